$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column G (shifts old G->H, old H->I)
$ws.Range("G:G").Insert()

# Copy formatting from the neighbouring header cell (F1) onto the new header cell (G1)
# so the new header gets the same bold/border/centered style used by the other headers.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Header text for the newly inserted "d=6" column
$ws.Range("G1").Value = "d=6"

# Data values for the newly inserted "d=6" column
$ws.Range("G2").Value = 97.82425019069851
$ws.Range("G3").Value = 97.89367936975371
$ws.Range("G4").Value = 97.83704822723142
$ws.Range("G5").Value = 97.80219795985222
$ws.Range("G6").Value = 97.81687888018982
